$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format so numeric-looking strings (e.g. "1.004", "0.06375")
# are not auto-converted to numbers by Excel's COM value-coercion.
$ws.Range('D2:E51').NumberFormat = "@"

$ws.Range('D2').Value = '25.949.64'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '1.643.28'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '215.52'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').Value = '0.5070'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = '1.004'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '0.2556'
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('D9').Value = '0.06379'
$ws.Range('E9').Value = '  -0.97%  '
$ws.Range('D10').Value = '19.49'
$ws.Range('E10').Value = '  -0.51%  '
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').Value = '1.657.94'
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('D13').Value = '4.282'
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('D14').Value = '0.5470'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').Value = '0.0₅7836'
$ws.Range('E15').Value = '  -1.29%  '
$ws.Range('D16').Value = '64.32'
$ws.Range('E16').Value = '  -0.06%  '
$ws.Range('D17').Value = '25.996.03'
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').Value = '197.49'
$ws.Range('E19').Value = '  -2.47%  '
$ws.Range('D20').Value = '4.451'
$ws.Range('E20').Value = '  +1.51%  '
$ws.Range('D21').Value = '9.949'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('D22').Value = '6.045'
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('D23').Value = '1.006'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = '1.898'
$ws.Range('E24').Value = '  +1.95%  '
$ws.Range('D25').Value = '140.81'
$ws.Range('E25').Value = '  -0.30%  '
$ws.Range('D26').Value = '0.1171'
$ws.Range('D27').Value = '6.896'
$ws.Range('E27').Value = '  +1.34%  '
$ws.Range('D28').Value = '15.71'
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').Value = '1.240'
$ws.Range('E29').Value = '  -0.47%  '
$ws.Range('D30').Value = '0.04958'
$ws.Range('E30').Value = '  +0.78%  '
$ws.Range('D31').Value = '3.262'
$ws.Range('E31').Value = '  -0.33%  '
$ws.Range('D32').Value = '3.187'
$ws.Range('E32').Value = '  -0.62%  '
$ws.Range('D33').Value = '1.541'
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('D34').Value = '2.369'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').Value = '0.8952'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  -1.61%  '
$ws.Range('D37').Value = '1.134.48'
$ws.Range('E37').Value = '  -1.58%  '
$ws.Range('D38').Value = '0.5439'
$ws.Range('E38').Value = '  -3.00%  '
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('D40').Value = '2.558'
$ws.Range('E40').Value = '  +0.17%  '
$ws.Range('D41').Value = '1.005'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '5.594'
$ws.Range('E42').Value = '  -1.98%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '0.8196'
$ws.Range('E43').Value = '  +1.57%  '
$ws.Range('E44').Value = '  +7.30%  '
$ws.Range('D45').Value = '99.57'
$ws.Range('E45').Value = '  -0.27%  '
$ws.Range('D46').Value = '1.777.58'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').Value = '0.4545'
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('D48').Value = '1.005'
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('D49').Value = '54.82'
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('D50').Value = '0.05080'
$ws.Range('E50').Value = '  +0.52%  '
$ws.Range('D51').Value = '1.006'
$ws.Range('E51').Value = '  +0.45%  '
